$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking text (e.g. '208.12') that must
# stay text, matching the workbook's existing inlineStr/shared-string cells.
# Force text format before assigning so Excel doesn't auto-coerce to Number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.118.84'
$ws.Range("E2").Value = '  +1.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.567.94'
$ws.Range("E3").Value = '  +1.95%  '
$ws.Range("E4").Value = '  +0.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.12'
$ws.Range("E5").Value = '  +1.23%  '
$ws.Range("E6").Value = '  +1.00%  '
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.09'
$ws.Range("E8").Value = '  +3.78%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.249'
$ws.Range("E9").Value = '  +1.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0588'
$ws.Range("E10").Value = '  +1.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0859'
$ws.Range("E11").Value = '  +0.60%  '
$ws.Range("E12").Value = '  +1.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.558.47'
$ws.Range("E13").Value = '  +1.27%  '
$ws.Range("E14").Value = '  +2.74%  '
$ws.Range("E15").Value = '  +2.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.099.30'
$ws.Range("E16").Value = '  +1.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.05'
$ws.Range("E17").Value = '  +1.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '219.46'
$ws.Range("E18").Value = '  +2.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0699'
$ws.Range("E19").Value = '  +2.64%  '
$ws.Range("E20").Value = '  +1.55%  '
$ws.Range("E22").Value = '  +1.91%  '
$ws.Range("E23").Value = '  +1.54%  '
$ws.Range("E24").Value = '  +1.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.26'
$ws.Range("E25").Value = '  +1.24%  '
$ws.Range("E26").Value = '  +0.82%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.00'
$ws.Range("E27").Value = '  +1.40%  '
$ws.Range("E28").Value = '  +0.34%  '
$ws.Range("E29").Value = '  +1.55%  '
$ws.Range("E30").Value = '  +2.98%  '
$ws.Range("E31").Value = '  +0.65%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.25'
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.454.06'
$ws.Range("E33").Value = '  +6.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.07'
$ws.Range("E34").Value = '  +4.79%  '
$ws.Range("E35").Value = '  +4.33%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.967'
$ws.Range("E36").Value = '  +1.05%  '
$ws.Range("E37").Value = '  +0.98%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0166'
$ws.Range("E38").Value = '  +0.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.525'
$ws.Range("E39").Value = '  +0.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.817'
$ws.Range("E40").Value = '  +1.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.76'
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("E42").Value = '  +0.32%  '
$ws.Range("E43").Value = '  +3.77%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.990'
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("E45").Value = '  +2.72%  '
$ws.Range("E46").Value = '  +1.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.705.79'
$ws.Range("E47").Value = '  +2.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.88'
$ws.Range("E48").Value = '  +3.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0525'
$ws.Range("E49").Value = '  +3.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0100'
$ws.Range("E50").Value = '  +2.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0966'
$ws.Range("E51").Value = '  +2.57%  '
